$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the "before" values of the columns that move (D, L, M, N, O, P, R, S)
# for every data row (2-41), then reassign them according to the row permutation
# described by the diff (rows were re-shuffled / re-dated while keeping the same
# set of record values).
$cols = @("D","L","M","N","O","P","R","S")
$snapshot = @{}
for ($r = 2; $r -le 41; $r++) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# Mapping: new row number -> source (old) row number
$mapping = @{
    2 = 24
    3 = 28
    4 = 30
    5 = 27
    6 = 9
    7 = 35
    8 = 6
    9 = 16
    10 = 36
    11 = 22
    12 = 17
    13 = 14
    14 = 29
    15 = 3
    16 = 2
    17 = 26
    18 = 20
    19 = 10
    20 = 21
    21 = 33
    22 = 37
    23 = 34
    24 = 15
    25 = 38
    26 = 40
    27 = 32
    28 = 23
    29 = 13
    30 = 4
    31 = 8
    32 = 25
    33 = 11
    34 = 41
    35 = 5
    36 = 31
    37 = 7
    38 = 18
    39 = 19
    40 = 12
    41 = 39
}

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $src = $snapshot[$oldRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value = $src[$c]
    }
}

$wb.Save()
